$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 ("Experimental") : Value column (B7) was empty -> now literal text "false".
# A plain Value assignment of "false" is auto-coerced to the Boolean FALSE by
# Excel's input parser, so build it as a text formula result and then paste
# the computed value back over itself (Paste Values) to "flatten" it into a
# genuine text cell without disturbing the existing cell style.
$cExperimental = $ws.Cells.Item(7, 2)
$cExperimental.Formula = "=""false"""
$cExperimental.Copy()
$cExperimental.PasteSpecial(-4163)  # xlPasteValues

# Row 8 ("Date") : refresh the publication date/time stamp.
$ws.Cells.Item(8, 2).Value = "2025-11-30T13:08:37+00:00"
